$p = $ppt.ActivePresentation

function Update-SusanSlide {
    param($Slide, $TitleSuffix)

    # "Subtitle 2" is the second shape on these title/closing slides.
    $shp = $Slide.Shapes.Item(2)
    $tr = $shp.TextFrame.TextRange

    # "Susan IBACH" -> "Susan IBACH <suffix>canada"
    $nameRange = $tr.Find("Susan IBACH")
    $nameRange.Text = "Susan IBACH" + $TitleSuffix + "canada"

    # Swap the hyperlinked old e-mail run for the new plain-text Amazon address.
    $emailRange = $tr.Find("SUSAN.IBACH@LIVE.COM")
    $emailRange.ActionSettings.Item(1).Hyperlink.Address = ""
    $emailRange.Text = "ibacsusa@amazon.com"

    # Drop the stray trailing space run that used to follow the hyperlink run.
    $newEmailRange = $tr.Find("ibacsusa@amazon.com")
    $tailStart = $newEmailRange.Start + $newEmailRange.Length
    if ($tailStart -le $tr.Length) {
        $tail = $tr.Characters($tailStart, 1)
        if ($tail.Text -eq " ") {
            $tail.Delete()
        }
    }
}

# Slide 1 (title slide): en-dash separated title, capitalized "Future".
$slide1Suffix = [string][char]0x20 + [char]0x2013 + " Amazon Future engineer program lead, "
Update-SusanSlide $p.Slides.Item(1) $slide1Suffix

# Slide 15 (closing "Questions?" slide): no dash, lower-case "future".
$slide15Suffix = " Amazon future engineer program lead, "
Update-SusanSlide $p.Slides.Item(15) $slide15Suffix
